{"js": "// Update the date line and every \"divided by\" answer cell in the table.\n// Each original value is unique in the document, so a simple\n// search + replace (exact match) for each pair is safe.\n\nconst replacements = [\n    { oldText: \"2024-08-09 Friday\", newText: \"2024-08-10 Saturday\" },\n    { oldText: \"830\u00f79=92, 2\",       newText: \"852\u00f72=426, 0\" },\n    { oldText: \"367\u00f77=52, 3\",       newText: \"525\u00f75=105, 0\" },\n    { oldText: \"517\u00f75=103, 2\",      newText: \"890\u00f76=148, 2\" },\n    { oldText: \"293\u00f75=58, 3\",       newText: \"486\u00f75=97, 1\" },\n    { oldText: \"220\u00f78=27, 4\",       newText: \"647\u00f78=80, 7\" },\n    { oldText: \"806\u00f78=100, 6\",      newText: \"898\u00f79=99, 7\" },\n    { oldText: \"692\u00f75=138, 2\",      newText: \"245\u00f75=49, 0\" },\n    { oldText: \"646\u00f79=71, 7\",       newText: \"330\u00f74=82, 2\" },\n    { oldText: \"258\u00f72=129, 0\",      newText: \"706\u00f75=141, 1\" },\n    { oldText: \"321\u00f76=53, 3\",       newText: \"586\u00f75=117, 1\" },\n    { oldText: \"269\u00f79=29, 8\",       newText: \"265\u00f79=29, 4\" },\n    { oldText: \"244\u00f79=27, 1\",       newText: \"439\u00f74=109, 3\" },\n    { oldText: \"748\u00f77=106, 6\",      newText: \"346\u00f78=43, 2\" },\n    { oldText: \"785\u00f73=261, 2\",      newText: \"605\u00f76=100, 5\" },\n    { oldText: \"179\u00f72=89, 1\",       newText: \"424\u00f72=212, 0\" },\n    { oldText: \"505\u00f79=56, 1\",       newText: \"261\u00f77=37, 2\" },\n    { oldText: \"836\u00f79=92, 8\",       newText: \"829\u00f75=165, 4\" },\n    { oldText: \"129\u00f72=64, 1\",       newText: \"200\u00f76=33, 2\" },\n    { oldText: \"735\u00f76=122, 3\",      newText: \"948\u00f77=135, 3\" },\n    { oldText: \"459\u00f78=57, 3\",       newText: \"809\u00f72=404, 1\" },\n    { oldText: \"646\u00f78=80, 6\",       newText: \"206\u00f74=51, 2\" },\n    { oldText: \"517\u00f73=172, 1\",      newText: \"812\u00f78=101, 4\" },\n    { oldText: \"295\u00f79=32, 7\",       newText: \"800\u00f74=200, 0\" },\n    { oldText: \"934\u00f75=186, 4\",      newText: \"433\u00f76=72, 1\" },\n    { oldText: \"186\u00f74=46, 2\",       newText: \"826\u00f75=165, 1\" }\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "# Update the date line and every \"divided by\" answer cell in the table.\n# Each original value is unique in the document, so a simple\n# Find/Replace (whole story, non-wildcard) for each pair is safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-08-09 Friday\"; New = \"2024-08-10 Saturday\" },\n    @{ Old = \"830\u00f79=92, 2\";       New = \"852\u00f72=426, 0\" },\n    @{ Old = \"367\u00f77=52, 3\";       New = \"525\u00f75=105, 0\" },\n    @{ Old = \"517\u00f75=103, 2\";      New = \"890\u00f76=148, 2\" },\n    @{ Old = \"293\u00f75=58, 3\";       New = \"486\u00f75=97, 1\" },\n    @{ Old = \"220\u00f78=27, 4\";       New = \"647\u00f78=80, 7\" },\n    @{ Old = \"806\u00f78=100, 6\";      New = \"898\u00f79=99, 7\" },\n    @{ Old = \"692\u00f75=138, 2\";      New = \"245\u00f75=49, 0\" },\n    @{ Old = \"646\u00f79=71, 7\";       New = \"330\u00f74=82, 2\" },\n    @{ Old = \"258\u00f72=129, 0\";      New = \"706\u00f75=141, 1\" },\n    @{ Old = \"321\u00f76=53, 3\";       New = \"586\u00f75=117, 1\" },\n    @{ Old = \"269\u00f79=29, 8\";       New = \"265\u00f79=29, 4\" },\n    @{ Old = \"244\u00f79=27, 1\";       New = \"439\u00f74=109, 3\" },\n    @{ Old = \"748\u00f77=106, 6\";      New = \"346\u00f78=43, 2\" },\n    @{ Old = \"785\u00f73=261, 2\";      New = \"605\u00f76=100, 5\" },\n    @{ Old = \"179\u00f72=89, 1\";       New = \"424\u00f72=212, 0\" },\n    @{ Old = \"505\u00f79=56, 1\";       New = \"261\u00f77=37, 2\" },\n    @{ Old = \"836\u00f79=92, 8\";       New = \"829\u00f75=165, 4\" },\n    @{ Old = \"129\u00f72=64, 1\";       New = \"200\u00f76=33, 2\" },\n    @{ Old = \"735\u00f76=122, 3\";      New = \"948\u00f77=135, 3\" },\n    @{ Old = \"459\u00f78=57, 3\";       New = \"809\u00f72=404, 1\" },\n    @{ Old = \"646\u00f78=80, 6\";       New = \"206\u00f74=51, 2\" },\n    @{ Old = \"517\u00f73=172, 1\";      New = \"812\u00f78=101, 4\" },\n    @{ Old = \"295\u00f79=32, 7\";       New = \"800\u00f74=200, 0\" },\n    @{ Old = \"934\u00f75=186, 4\";      New = \"433\u00f76=72, 1\" },\n    @{ Old = \"186\u00f74=46, 2\";       New = \"826\u00f75=165, 1\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
